$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, avoiding Excel auto-converting
# numeric-looking strings (e.g. "226.48") into real numbers, and without
# leaving behind a residual "Text" number-format style on the cell.
function Set-TextValue($cell, $text) {
    $ws.Range($cell).Value = "'$text"
    $ws.Range($cell).Style = "Normal"
}

$ws.Range('D2').Value = '34.049.35'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.787.51'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '226.48'
$ws.Range('E5').Value = '  +1.99%  '
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('E7').Value = '  +0.08%  '
Set-TextValue 'D8' '32.16'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  +3.68%  '
Set-TextValue 'D10' '0.0682'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').Value = '2.045.19'
$ws.Range('E12').Value = '  -0.05%  '
Set-TextValue 'D13' '11.28'
$ws.Range('E13').Value = '  +3.16%  '
$ws.Range('D14').Value = '1.791.74'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('D15').Value = '34.017.90'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('E16').Value = '  -1.17%  '
Set-TextValue 'D18' '67.71'
$ws.Range('E18').Value = '  -0.48%  '
Set-TextValue 'D19' '242.58'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = '0.0₃0772'
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('E24').Value = '  -2.80%  '
Set-TextValue 'D25' '161.90'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('E26').Value = '  +1.13%  '
Set-TextValue 'D27' '16.19'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('E32').Value = '  -0.99%  '
$ws.Range('E33').Value = '  +2.76%  '
Set-TextValue 'D34' '1.83'
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('D35').Value = '1.397.40'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D38' '0.0188'
$ws.Range('E38').Value = '  +1.51%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D39' '2.34'
$ws.Range('E39').Value = '  +7.95%  '
Set-TextValue 'D40' '79.92'
$ws.Range('E40').Value = '  +0.36%  '
Set-TextValue 'D41' '2.36'
$ws.Range('E41').Value = '  +0.05%  '
Set-TextValue 'D42' '0.918'
$ws.Range('E42').Value = '  -0.19%  '
Set-TextValue 'D43' '13.69'
$ws.Range('E43').Value = '  +14.03%  '
$ws.Range('E44').Value = '  -1.28%  '
$ws.Range('E45').Value = '  +8.16%  '
Set-TextValue 'D47' '6.06'
$ws.Range('E47').Value = '  +2.46%  '
$ws.Range('E48').Value = '  +2.83%  '
Set-TextValue 'D49' '107.59'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').Value = '1.946.03'
$ws.Range('E50').Value = '  -0.29%  '
